# Lucky Ocean review update:
#  - New title/headline (H1, the bold restatement near the end, and the meta
#    description).
#  - Reworked "What we like" bullet list (4 bullets, reordered/edited).
#  - Reworked "What we don't like" bullet list (2 bullets, edited).
#
# We use Range.InsertXML on each target paragraph's whole Range (which
# includes its <w:pPr>) so the exact paragraph formatting (list style,
# spacing, indentation, run formatting) is preserved while we fully control
# the resulting run structure (including the placeholder empty <w:r/> that
# precedes the text run in the original document).

$d = $word.ActiveDocument

function Set-ParagraphInnerXml($paraIndex, $innerXml) {
    $p = $d.Paragraphs.Item($paraIndex)
    $range = $p.Range
    $pkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $innerXml + '</w:document></pkg:xmlData></pkg:part></pkg:package>'
    $range.InsertXML($pkg)
}

function Heading1Para($text) {
    return '<w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>' + $text + '</w:t></w:r></w:p>'
}

function BulletPara($text) {
    return '<w:p><w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr><w:r/><w:r><w:t>' + $text + '</w:t></w:r></w:p>'
}

function BoldPara($text) {
    return '<w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>' + $text + '</w:t></w:r></w:p>'
}

function ItalicPara($text) {
    return '<w:p><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>' + $text + '</w:t></w:r></w:p>'
}

$newTitle = "Play Lucky Ocean and Win Big - Free Online Slot Game"

# 1. Main H1 title (paragraph 1).
$p1 = Heading1Para $newTitle
$body1 = '<w:body>' + $p1 + '</w:body>'
Set-ParagraphInnerXml 1 $body1

# 2. "What we like" bullets (paragraphs 39-42), rewritten in their final
#    order: mechanics, intuitive, interesting winnings, growing jackpot.
$p39 = BulletPara "Unique instant lottery mechanics"
$body39 = '<w:body>' + $p39 + '</w:body>'
Set-ParagraphInnerXml 39 $body39

$p40 = BulletPara "Intuitive and enjoyable for new players"
$body40 = '<w:body>' + $p40 + '</w:body>'
Set-ParagraphInnerXml 40 $body40

$p41 = BulletPara "Interesting winnings depend on player's choices"
$body41 = '<w:body>' + $p41 + '</w:body>'
Set-ParagraphInnerXml 41 $body41

$p42 = BulletPara "Growing jackpot as player unlocks more pearls"
$body42 = '<w:body>' + $p42 + '</w:body>'
Set-ParagraphInnerXml 42 $body42

# 3. "What we don't like" bullets (paragraphs 44-45).
$p44 = BulletPara "Generic sound that could better align with the theme"
$body44 = '<w:body>' + $p44 + '</w:body>'
Set-ParagraphInnerXml 44 $body44

$p45 = BulletPara "Losing all collected and unredeemed loot if Porcupine Fish is found"
$body45 = '<w:body>' + $p45 + '</w:body>'
Set-ParagraphInnerXml 45 $body45

# 4. Bold restatement of the title (paragraph 46).
$p46 = BoldPara $newTitle
$body46 = '<w:body>' + $p46 + '</w:body>'
Set-ParagraphInnerXml 46 $body46

# 5. Italic meta description (paragraph 47).
$newMeta = "Play Lucky Ocean, a unique online slot game with growing jackpot. Win big and play for free!"
$p47 = ItalicPara $newMeta
$body47 = '<w:body>' + $p47 + '</w:body>'
Set-ParagraphInnerXml 47 $body47
